$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Centric Tops Out at 58-01 Queens Boulevard in Woodside, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/09/centric-tops-out-at-58-01-queens-boulevard-in-woodside-queens.html"
$ws.Range("C2").Value = 'Construction has topped out on Centric, a 12-story residential building at 58-01 Queens Boulevard in <a href="https://newyorkyimby.com/neighborhoods/woodside">Woodside</a>, Queens. Designed by Tang Studio Architect and developed by New Empire Corporation, the structure will yield 131 condominium units in one- to three-bedroom layouts, as well as a collection of townhomes. The property is bounded by 44th Avenue to the north, Queens Boulevard to the south, and 58th Street to the west.'
$ws.Range("D2").Value = "2025-09-14T11:30:14+00:00"
$ws.Range("E2").Value = "Sun, 14 Sep 2025 11:30:14 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Woodside"
$ws.Range("H2").Style = "Normal"

$ws.Range("A3").Value = "New Permits Filed for 27-05 27th Street in Astoria, Queens"
$ws.Range("B3").Value = "https://newyorkyimby.com/2025/09/new-permits-filed-for-27-05-27th-street-in-astoria-queens.html"
$ws.Range("C3").Value = "New permits have been filed for a four-story residential building at 27-05 27th Street in <a href=`"https://newyorkyimby.com/neighborhoods/Astoria`">Astoria</a>, Queens. Located at the intersection of Newtown Avenue and 27th Street, the lot is three blocks west of the 30th Avenue subway station, served by the N and W trains. Joe D'Amico of Citia Nova, LLC is listed as the owner behind the applications."
$ws.Range("D3").Value = "2025-09-14T10:30:36+00:00"
$ws.Range("E3").Value = "Sun, 14 Sep 2025 10:30:36 +0000"
$ws.Range("F3").Value = "YIMBY"
$ws.Range("G3").Value = "YIMBY - Astoria"
$ws.Range("H3").Style = "Normal"
